# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 365   # was 362 - 南宁·原神x星铁x绝区零同人ONLY3.0
$wsExpo.Range("F3").Value = 68    # was 66  - 南宁·小野爷爷&娃展2.0
$wsExpo.Range("F4").Value = 275   # was 274 - 南宁·布谷鸟动漫展5th
$wsExpo.Range("F5").Value = 4101  # was 4088 - 南宁·2024良牙动漫秋季盛典（秋典）
$wsExpo.Range("F7").Value = 450   # was 449 - 南宁·万圣漫控嘉年华10

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 365    # was 362 - 南宁·原神x星铁x绝区零同人ONLY3.0
$wsAll.Range("F3").Value = 68     # was 66  - 南宁·小野爷爷&娃展2.0
$wsAll.Range("F4").Value = 275    # was 274 - 南宁·布谷鸟动漫展5th
$wsAll.Range("F5").Value = 4101   # was 4088 - 南宁·2024良牙动漫秋季盛典（秋典）
$wsAll.Range("F9").Value = 450    # was 449 - 南宁·万圣漫控嘉年华10
